# Auto-generated: update cryptos list (Price / Volume(1h)) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.950.77'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '1.556.66'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.05'
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.06'
$ws.Range("E8").Value = '  +2.90%  '
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("E10").Value = '  +0.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0857'
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").Value = '1.778.08'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").Value = '1.555.56'
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.75'
$ws.Range("E14").Value = '  +1.20%  '
$ws.Range("E15").Value = '  +1.68%  '
$ws.Range("D16").Value = '26.954.21'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.79'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '217.86'
$ws.Range("E18").Value = '  +1.37%  '
$ws.Range("E19").Value = '  +1.73%  '
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("E22").Value = '  +1.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.24'
$ws.Range("E23").Value = '  +0.55%  '
$ws.Range("E24").Value = '  +0.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.17'
$ws.Range("E25").Value = '  +1.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.65'
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.94'
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("E30").Value = '  +2.34%  '
$ws.Range("E31").Value = '  -0.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.23'
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("D33").Value = '1.424.59'
$ws.Range("E33").Value = '  +3.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.08'
$ws.Range("E34").Value = '  +4.25%  '
$ws.Range("E35").Value = '  +3.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.976'
$ws.Range("E36").Value = '  +2.10%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.522'
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("E42").Value = '  +1.03%  '
$ws.Range("E43").Value = '  +3.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.984'
$ws.Range("E44").Value = '  -0.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.70'
$ws.Range("E45").Value = '  +1.61%  '
$ws.Range("E46").Value = '  +0.66%  '
$ws.Range("D47").Value = '1.692.32'
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.41'
$ws.Range("E48").Value = '  +2.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0521'
$ws.Range("E49").Value = '  +3.11%  '
$ws.Range("E50").Value = '  +3.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0960'
$ws.Range("E51").Value = '  +1.25%  '
